# Edit script: applies the commit's changes to
# relatorio_comparacao_rest_grpc.docx via the Word COM object model.
#
# Summary of the edit:
#   1. Adds a new "Nome: ..." paragraph right after the title.
#   2. Rewrites the introductory paragraph and appends hardware info
#      (MacBook Air M1, RAM, SSD) as manual line breaks in the same run.
#   3. Updates the lead-in sentence before the benchmark results.
#   4. Replaces the REST/gRPC bullet list + closing paragraph with a
#      4-column results table (Operacao/Metrica/REST/gRPC).

$d = $word.ActiveDocument

# --- sanity check on the pristine document shape we expect to edit ---
$expectedP2 = "O presente relatório tem como objetivo comparar o uso das abordagens REST e gRPC em termos de latência e throughput, com base em dados obtidos a partir de benchmarks executados no contexto da disciplina."
if ($d.Paragraphs(2).Range.Text.TrimEnd([char]13) -ne $expectedP2) {
    throw "Unexpected paragraph 2 content; aborting to avoid corrupting the document."
}

# 1) Insert the new "Nome:" paragraph before the intro paragraph (old paragraph 2)
$d.Paragraphs(2).Range.InsertParagraphBefore()
$d.Paragraphs(2).Range.Text = "Nome: Eric Rodrigues Diniz"

# 2) Replace the intro paragraph text (now paragraph 3) with the new
#    multi-line text. [char]11 is a manual line break (<w:br/>) inside
#    the same run, matching the target markup.
$nl = [char]11
$introText = "O relatório a seguir tem como objetivo comparar o uso das abordagens REST e gRPC em termos de latência e throughput, com base em dados obtidos a partir de benchmarks executados diretamente no meu computador pessoal:" + $nl + $nl + "MacBook Air M1 (2020)" + $nl + "8 GB RAM" + $nl + "256 GB de armazenamento SSD"
$d.Paragraphs(3).Range.Text = $introText

# 3) Update the "A seguir, ..." lead-in sentence before the benchmark results
$d.Content.Find.Execute("A seguir, são apresentados alguns resultados coletados nos testes realizados com gRPC:", $true, $false, $false, $false, $false, $true, 1, $false, "A seguir, são apresentados os resultados coletados nos testes realizados com REST e gRPC:", 2) | Out-Null

# 4) Remove the old bullet list (LIST/CREATE/STATS) plus the closing
#    paragraph, then insert the new results table in their place.
#    After steps 1-3 the document layout is:
#      1 Heading1 / 2 Nome / 3 intro / 4 REST heading / 5 REST text /
#      6 gRPC heading / 7 gRPC text / 8 Resultados heading /
#      9 "A seguir..." / 10-12 bullets / 13 "Esses valores..." / 14 Conclusao
$expectedP13 = "Esses valores demonstram que o gRPC é capaz de responder rapidamente mesmo sob concorrência elevada, oferecendo desempenho superior ao esperado em implementações REST tradicionais."
if ($d.Paragraphs(13).Range.Text.TrimEnd([char]13) -ne $expectedP13) {
    throw "Unexpected paragraph 13 content; aborting to avoid corrupting the document."
}

$delStart = $d.Paragraphs(10).Range.Start
$delEnd = $d.Paragraphs(14).Range.Start
$delRange = $d.Range($delStart, $delEnd)
$delRange.Delete()

$insPoint = $d.Range($delStart, $delStart)
$tblXml = '<w:tbl><w:tblPr><w:tblW w:type="auto" w:w="0"/><w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/></w:tblPr><w:tblGrid><w:gridCol w:w="2160"/><w:gridCol w:w="2160"/><w:gridCol w:w="2160"/><w:gridCol w:w="2160"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Operação</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Métrica</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>REST (medido)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>gRPC (medido)</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>LIST</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>p50 (ms)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~10</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~10.9</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>LIST</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>p95/p97.5 (ms)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~21 (p97.5)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~12.9 (p95)</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>LIST</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Throughput (rps)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~500</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~1831</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>CREATE</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>p50 (ms)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~15</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~5.6</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>CREATE</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>p95/p97.5 (ms)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~66 (p97.5)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~9.3 (p95)</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>CREATE</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Throughput (rps)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~500</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~1663</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>STATS</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>p50 (ms)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~22.5</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>STATS</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>p95/p97.5 (ms)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~7 (p97.5)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~24.7 (p95)</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>STATS</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Throughput (rps)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~500</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>~2303</w:t></w:r></w:p></w:tc></w:tr></w:tbl>'
$pkgXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $tblXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint.InsertXML($pkgXml) | Out-Null

Write-Output "edit applied"
